$d = $word.ActiveDocument

# 1. Remove the existing _GoBack bookmark - it will be reinserted at its new location
try {
  $d.Bookmarks("_GoBack").Delete()
} catch {
}

# 2. Mint a new list numbering definition (will become numId=5 / a fresh abstractNum)
#    by applying default numbering to the last (empty) paragraph; we immediately
#    overwrite that paragraph's content below, but the numbering definitions
#    it creates in numbering.xml persist for reuse via explicit numId references.
$mintPara = $d.Paragraphs.Last
$mintPara.Range.ListFormat.ApplyNumberDefault()

# 3. Replace the content of that same (still last, still empty of real text) paragraph
#    with the full block of new content described in the change: a section-break
#    paragraph, the "Analisis Encriptado2.txt" heading, and the supporting bullet /
#    numbered paragraphs, ending with the relocated _GoBack bookmark.
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range

$xml = @'
    <w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:ind w:left="1080"/>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:sectPr>
          <w:pgSz w:w="12240" w:h="15840"/>
          <w:pgMar w:top="1417" w:right="1701" w:bottom="1417" w:left="1701" w:header="708" w:footer="708" w:gutter="0"/>
          <w:cols w:space="708"/>
          <w:docGrid w:linePitch="360"/>
        </w:sectPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>Análisis Encriptado2.txt</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve">Revisando el texto encriptado, mediante la observación me pude dar cuenta que también se </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>repetía</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> el mismo patrón de la Z, esta vez la Z estaba presente en cada tripleta, en algunas 2 veces repetidas y en otras solo una vez, pero nunca se repetía 3 veces</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="4"/>
        </w:numPr>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve">Al igual que el Encriptado1.txt asumí que la Z me estaba dando información acerca </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>del</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> XOR aplicado (0x5A = 90 = Z), pero a diferencia del primer encriptado la Z también me estaba indicando que letra o que pareja se estaba repitiendo, porque hay que recordar que en este problema se utilizó el método de compresión LZ78, entonces mediante unas pruebas manuales pude corroborar que:</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="5"/>
        </w:numPr>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve">Las tripletas me indican que </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>cuando hay dos Z seguidas el tercer carácter no se repite, en cambio cuando la Z tiene otra letra distinta a ella en la segunda posición significa que esa letra ya estaba repetida, por ende, siguiendo la lógica de LZ78 esa letra repetida equivale a un índice a el cual le vamos a agregar una pareja (otra letra).</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="5"/>
        </w:numPr>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve">También pude corroborar que el primer índice repetido es igual a J, siguiendo la lógica de LZ78 el primer índice repetido es el “a”, </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>entonces</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> a ese índice repetido le agrego “n” = nuevo </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>índice</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> “an”.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
      <w:pPr>
        <w:pStyle w:val="Prrafodelista"/>
        <w:numPr>
          <w:ilvl w:val="0"/>
          <w:numId w:val="5"/>
        </w:numPr>
        <w:rPr>
          <w:b/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>Sabiendo eso y comprobando que coincide con el archivo Encriptado2.txt, surge una pregunta</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t>Qué</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:lang w:val="es-ES"/>
        </w:rPr>
        <w:t xml:space="preserve"> operación oculta se aplicó para  que el índice a sea igual a J?</w:t>
      </w:r>
    </w:p>
'@

$r.InsertXML($xml)

Write-Output ("Paragraphs count: " + $d.Paragraphs.Count)
